# Fixed Bento 80 Test scripts: append/adjust Cypher "order by ... limit" clauses
# for the 3 saved queries on the "startup" sheet, grow the wrapped-text row
# heights to fit the extra line, and leave the selection on the last query cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current query text for the three query cells (Value2 avoids the
# property-descriptor quirk of a bare .Value read in this host).
$b2 = $ws.Range("B2").Value2
$b3 = $ws.Range("B3").Value2
$b4 = $ws.Range("B4").Value2

# CasesTab query: add ordering + limit after the last returned column.
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# SamplesTab query: add ordering + limit after the last returned column.
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# FilesTab query: replace the existing lowercase "order by" with the new
# capitalized clause that also adds ASC + LIMIT 100.
$b4New = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value2 = $b4New

# The extra line in each query needs a taller wrapped row.
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360

# Leave the selection on the Files query cell, as in the saved workbook.
$ws.Range("B4").Select()
